# Add placeholder replacement in headers and footers of Word documents.
#
# The document's default header contains merge-style placeholder tokens
# (e.g. "{PlaceHolder2}", "{PlaceHolder3}") left over from a template.
# Replace each placeholder with its real value, searching every header
# and footer story (default / first-page / even-page) of every section
# so the substitution is applied uniformly wherever it occurs.

$d = $word.ActiveDocument

# Map of placeholder token -> replacement value.
$replacements = @{
    "{PlaceHolder2}" = "90.54"
    "{PlaceHolder3}" = "3245789085"
}

foreach ($sec in $d.Sections) {
    $stories = @()
    $stories += $sec.Headers
    $stories += $sec.Footers

    foreach ($story in $stories) {
        if (-not $story.Exists) { continue }

        foreach ($key in $replacements.Keys) {
            $story.Range.Find.Execute(
                $key, $true, $false, $false, $false, $false,
                $true, 1, $false, $replacements[$key], 2
            )
        }
    }
}

Write-Host "Placeholder replacement complete."
